$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(19, 8).Value = 3197.75
$ws.Cells.Item(19, 9).Value = 894.5
$ws.Cells.Item(19, 10).Value = 5501
$ws.Cells.Item(19, 11).Value = 894.5
$ws.Cells.Item(19, 12).Value = 5501
$ws.Cells.Item(19, 13).Value = -719.5
$ws.Cells.Item(19, 14).Value = -5851
$ws.Cells.Item(51, 8).Value = 13892222
$ws.Cells.Item(51, 9).Value = 2499.5
$ws.Cells.Item(51, 11).Value = 2499.5
$ws.Cells.Item(51, 13).Value = -2015.5
$ws.Cells.Item(113, 8).Value = 6680
$ws.Cells.Item(113, 9).Value = 12400
$ws.Cells.Item(113, 10).Value = 5250
$ws.Cells.Item(113, 11).Value = 12400
$ws.Cells.Item(113, 12).Value = 5250
$ws.Cells.Item(113, 13).Value = -9146
$ws.Cells.Item(113, 14).Value = -11758
$ws.Cells.Item(138, 8).Value = 2394.6858
$ws.Cells.Item(138, 10).Value = 2910.25
$ws.Cells.Item(138, 12).Value = 8730.75
$ws.Cells.Item(138, 14).Value = -19010.75

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(28, 8).Value = 12822.1
$ws.Cells.Item(28, 9).Value = 6740.4287
$ws.Cells.Item(28, 11).Value = 6740.4287
$ws.Cells.Item(28, 13).Value = -6548.4287
$ws.Cells.Item(74, 8).Value = 1846.7059
$ws.Cells.Item(74, 9).Value = 1649.6875
$ws.Cells.Item(74, 10).Value = 4999
$ws.Cells.Item(74, 11).Value = 1649.6875
$ws.Cells.Item(74, 12).Value = 4999
$ws.Cells.Item(74, 13).Value = -775.6875
$ws.Cells.Item(74, 14).Value = -6747
$ws.Cells.Item(77, 8).Value = 1846.7059
$ws.Cells.Item(77, 9).Value = 1649.6875
$ws.Cells.Item(77, 10).Value = 4999
$ws.Cells.Item(77, 11).Value = 8248.4375
$ws.Cells.Item(77, 12).Value = 24995
$ws.Cells.Item(77, 13).Value = -3880.4375
$ws.Cells.Item(77, 14).Value = -33731
$ws.Cells.Item(99, 8).Value = 12822.1
$ws.Cells.Item(99, 9).Value = 6740.4287
$ws.Cells.Item(99, 11).Value = 6740.4287
$ws.Cells.Item(99, 13).Value = -3745.4287
$ws.Cells.Item(119, 8).Value = 58999
$ws.Cells.Item(119, 10).Value = 58999
$ws.Cells.Item(119, 12).Value = 58999
$ws.Cells.Item(119, 14).Value = -68675
$ws.Cells.Item(132, 8).Value = 1958.7073
$ws.Cells.Item(132, 9).Value = 1434.5927
$ws.Cells.Item(132, 11).Value = 4303.7781
$ws.Cells.Item(132, 13).Value = -1773.7781

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 8724.679
$ws.Cells.Item(20, 9).Value = 7826.5884
$ws.Cells.Item(20, 11).Value = 7826.5884
$ws.Cells.Item(20, 13).Value = -7579.5884
$ws.Cells.Item(98, 8).Value = 30542
$ws.Cells.Item(98, 10).Value = 30542
$ws.Cells.Item(98, 12).Value = 30542
$ws.Cells.Item(98, 14).Value = -36532
$ws.Cells.Item(99, 8).Value = 2036.625
$ws.Cells.Item(99, 9).Value = 1315.0769
$ws.Cells.Item(99, 10).Value = 5163.3335
$ws.Cells.Item(99, 11).Value = 1315.0769
$ws.Cells.Item(99, 12).Value = 5163.3335
$ws.Cells.Item(99, 13).Value = 182.9231
$ws.Cells.Item(99, 14).Value = -8159.3335
$ws.Cells.Item(134, 8).Value = 1348.25
$ws.Cells.Item(134, 9).Value = 1218
$ws.Cells.Item(134, 11).Value = 3654
$ws.Cells.Item(134, 13).Value = -1119

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 1649.05
$ws.Cells.Item(16, 9).Value = 1555.8572
$ws.Cells.Item(16, 10).Value = 1866.5
$ws.Cells.Item(16, 11).Value = 1555.8572
$ws.Cells.Item(16, 12).Value = 1866.5
$ws.Cells.Item(16, 13).Value = -1268.8572
$ws.Cells.Item(16, 14).Value = -2440.5
$ws.Cells.Item(31, 8).Value = 2405.6365
$ws.Cells.Item(31, 10).Value = 5000
$ws.Cells.Item(31, 12).Value = 5000
$ws.Cells.Item(31, 14).Value = -5590
$ws.Cells.Item(34, 8).Value = 2405.6365
$ws.Cells.Item(34, 10).Value = 5000
$ws.Cells.Item(34, 12).Value = 5000
$ws.Cells.Item(34, 14).Value = -5404
$ws.Cells.Item(86, 8).Value = 8892.143
$ws.Cells.Item(86, 9).Value = 8408.182000000001
$ws.Cells.Item(86, 11).Value = 8408.182000000001
$ws.Cells.Item(86, 13).Value = -7285.182000000001
$ws.Cells.Item(89, 8).Value = 8892.143
$ws.Cells.Item(89, 9).Value = 8408.182000000001
$ws.Cells.Item(89, 11).Value = 42040.91
$ws.Cells.Item(89, 13).Value = -36424.91
$ws.Cells.Item(105, 8).Value = 1854.2307
$ws.Cells.Item(105, 9).Value = 1840
$ws.Cells.Item(105, 10).Value = 1870.8334
$ws.Cells.Item(105, 11).Value = 1840
$ws.Cells.Item(105, 12).Value = 1870.8334
$ws.Cells.Item(105, 13).Value = -93
$ws.Cells.Item(105, 14).Value = -5364.8334
$ws.Cells.Item(113, 8).Value = 1649.05
$ws.Cells.Item(113, 9).Value = 1555.8572
$ws.Cells.Item(113, 10).Value = 1866.5
$ws.Cells.Item(113, 11).Value = 1555.8572
$ws.Cells.Item(113, 12).Value = 1866.5
$ws.Cells.Item(113, 13).Value = 614.1428000000001
$ws.Cells.Item(113, 14).Value = -6206.5
$ws.Cells.Item(122, 8).Value = 5482.727
$ws.Cells.Item(122, 9).Value = 5231
$ws.Cells.Item(122, 11).Value = 15693
$ws.Cells.Item(122, 13).Value = -13243
$ws.Cells.Item(132, 8).Value = 3999.5
$ws.Cells.Item(132, 9).Value = 3999.5
$ws.Cells.Item(132, 11).Value = 11998.5
$ws.Cells.Item(132, 13).Value = -9468.5
$ws.Cells.Item(134, 8).Value = 6352.125
$ws.Cells.Item(134, 9).Value = 6352.125
$ws.Cells.Item(134, 11).Value = 19056.375
$ws.Cells.Item(134, 13).Value = -16521.375

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(23, 8).Value = 265
$ws.Cells.Item(23, 10).Value = 306.25
$ws.Cells.Item(23, 12).Value = 918.75
$ws.Cells.Item(23, 14).Value = -1388.75
$ws.Cells.Item(34, 8).Value = 1281.0714
$ws.Cells.Item(34, 9).Value = 111.42857
$ws.Cells.Item(34, 10).Value = 2450.7144
$ws.Cells.Item(34, 11).Value = 334.28571
$ws.Cells.Item(34, 12).Value = 7352.1432
$ws.Cells.Item(34, 13).Value = -250.28571
$ws.Cells.Item(34, 14).Value = -7520.1432
$ws.Cells.Item(39, 8).Value = 4297.1
$ws.Cells.Item(39, 9).Value = 1748.5
$ws.Cells.Item(39, 10).Value = 4934.25
$ws.Cells.Item(39, 11).Value = 5245.5
$ws.Cells.Item(39, 12).Value = 14802.75
$ws.Cells.Item(39, 13).Value = -4951.5
$ws.Cells.Item(39, 14).Value = -15390.75
$ws.Cells.Item(55, 8).Value = 7814395.5
$ws.Cells.Item(55, 10).Value = 12502760
$ws.Cells.Item(55, 12).Value = 37508280
$ws.Cells.Item(55, 14).Value = -37508634

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(64, 8).Value = 23271
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 13).Value = $null
$ws.Cells.Item(67, 8).Value = 23271
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 13).Value = $null
$ws.Cells.Item(70, 8).Value = 5004.5
$ws.Cells.Item(70, 9).Value = 4500
$ws.Cells.Item(70, 11).Value = 4500
$ws.Cells.Item(70, 13).Value = -4230
$ws.Cells.Item(73, 8).Value = 5004.5
$ws.Cells.Item(73, 9).Value = 4500
$ws.Cells.Item(73, 11).Value = 4500
$ws.Cells.Item(73, 13).Value = -3564
$ws.Cells.Item(132, 8).Value = 2513.111
$ws.Cells.Item(132, 9).Value = 2284.4
$ws.Cells.Item(132, 10).Value = 3656.6667
$ws.Cells.Item(132, 11).Value = 6853.200000000001
$ws.Cells.Item(132, 12).Value = 10970.0001
$ws.Cells.Item(132, 13).Value = -4323.200000000001
$ws.Cells.Item(132, 14).Value = -16030.0001

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value = 1396.3334
$ws.Cells.Item(22, 9).Value = 1609.6666
$ws.Cells.Item(22, 11).Value = 1609.6666
$ws.Cells.Item(22, 13).Value = -1314.6666
$ws.Cells.Item(27, 8).Value = 1396.3334
$ws.Cells.Item(27, 9).Value = 1609.6666
$ws.Cells.Item(27, 11).Value = 1609.6666
$ws.Cells.Item(27, 13).Value = -1502.6666
$ws.Cells.Item(99, 8).Value = 25000
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 14).Value = $null
$ws.Cells.Item(100, 8).Value = 22146.15
$ws.Cells.Item(100, 9).Value = 6204.1665
$ws.Cells.Item(100, 11).Value = 6204.1665
$ws.Cells.Item(100, 13).Value = -5663.1665
$ws.Cells.Item(132, 8).Value = 3754500
$ws.Cells.Item(132, 9).Value = 7500000
$ws.Cells.Item(132, 10).Value = 9000
$ws.Cells.Item(132, 11).Value = 22500000
$ws.Cells.Item(132, 12).Value = 27000
$ws.Cells.Item(132, 13).Value = -22497470
$ws.Cells.Item(132, 14).Value = -32060
$ws.Cells.Item(136, 8).Value = 3076.6365
$ws.Cells.Item(136, 9).Value = 2432.8333
$ws.Cells.Item(136, 11).Value = 7298.499899999999
$ws.Cells.Item(136, 13).Value = -4748.499899999999

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 3526
$ws.Cells.Item(81, 9).Value = 3470.2354
$ws.Cells.Item(81, 10).Value = 4000
$ws.Cells.Item(81, 11).Value = 6940.4708
$ws.Cells.Item(81, 12).Value = 8000
$ws.Cells.Item(81, 13).Value = -5879.4708
$ws.Cells.Item(81, 14).Value = -10122
$ws.Cells.Item(84, 8).Value = 3526
$ws.Cells.Item(84, 9).Value = 3470.2354
$ws.Cells.Item(84, 10).Value = 4000
$ws.Cells.Item(84, 11).Value = 34702.354
$ws.Cells.Item(84, 12).Value = 40000
$ws.Cells.Item(84, 13).Value = -29398.354
$ws.Cells.Item(84, 14).Value = -50608
$ws.Cells.Item(100, 8).Value = 2783.6
$ws.Cells.Item(100, 9).Value = 2842.6155
$ws.Cells.Item(100, 11).Value = 5685.231
$ws.Cells.Item(100, 13).Value = -5144.231
$ws.Cells.Item(119, 8).Value = 24898
$ws.Cells.Item(119, 10).Value = 24898
$ws.Cells.Item(119, 12).Value = 24898
$ws.Cells.Item(119, 14).Value = -34574
$ws.Cells.Item(132, 8).Value = 100000
$ws.Cells.Item(132, 9).Value = 100000
$ws.Cells.Item(132, 11).Value = 300000
$ws.Cells.Item(132, 13).Value = -297470
$ws.Cells.Item(136, 8).Value = 265.83334
$ws.Cells.Item(136, 9).Value = 265.83334
$ws.Cells.Item(136, 11).Value = 797.5000200000001
$ws.Cells.Item(136, 13).Value = 1752.49998
